# Zestaponi municipality area sheet: rename the sheet, drop the old
# "census" subtitle/spacer rows, and collapse the 1989/2002/2014 table
# down to just the most recent (2014) figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the blank spacer row that used to sit above the year header
# (row 4/5/6 shift up to become row 3/4/5).
$ws.Rows("3").Delete()

# Keep only the 2014 column; the 1989/2002 columns go away (old column D
# shifts left into column B).
$ws.Columns("B:C").Delete()

# The old "(მოსახლეობის აღწერის შედეგებით)" subtitle under the title is
# removed, leaving row 2 blank.
$ws.Range("A2").Clear()

# Give the sheet its real name instead of the placeholder "1".
$ws.Name = "ზესტაფონი"

$ws.Range("A2").Select()
